# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) holds a per-row statistic (s_val) that is
# recomputed from the option-chain data for each row (strike count /
# distribution bucket). Column G1 already carries the "K" header, so
# here we just recalculate and re-write the numeric K values for every
# data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly computed K (s_val) numbers keyed by row number, as produced by
# the regenerated std/mean calculation pass.
$sVals = @{
    2  = 1;  3  = 2;  4  = 0;  5  = 0;  6  = 2;  7  = 2;  8  = 2;  9  = 1;  10 = 2;
    11 = 2;  12 = 0;  13 = 1;  14 = 1;  15 = 2;  16 = 1;  18 = 1;  19 = 1;  20 = 1;
    21 = 0;  22 = 3;  23 = 0;  24 = 2;  25 = 2;  26 = 1;  27 = 0;  28 = 1;  29 = 0;
    30 = 0;  31 = 1;  32 = 1;  33 = 0;  34 = 1;  35 = 2;  36 = 3;  37 = 2;  38 = 2;
    39 = 2;  40 = 1;  41 = 1;  42 = 1;  43 = 1;  44 = 1;  45 = 1;  46 = 2;  47 = 1;
    49 = 3;  50 = 3;  51 = 3;  52 = 1;  53 = 2;  54 = 0;  55 = 0;  56 = 1;  57 = 0;
    59 = 3;  60 = 1;  61 = 0;  62 = 1;  63 = 2;  64 = 1
}

foreach ($row in $sVals.Keys) {
    $ws.Range("G$row").Value = $sVals[$row]
}
